$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date value from 45192 (2023-09-23) to
# 45202 (2023-10-03) for every existing data row (rows 2-23).
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# The previously-last row (23) picks up an explicit row height once the
# sheet is touched/extended - force it so it matches the other data rows.
$ws.Rows.Item(23).RowHeight = 15

# Append the new case as row 24.
$ws.Cells.Item(24, 1).Value = "A 47107-2023"

$ws.Cells.Item(24, 2).Value = 45196
$ws.Cells.Item(24, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(24, 3).Value = 45202
$ws.Cells.Item(24, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(24, 4).Value = "OKÄNT"
$ws.Cells.Item(24, 5).Value = "OKÄNT"

$ws.Cells.Item(24, 7).Value = 1.5
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0

# Column R mirrors the other rows: an empty, wrap-text-styled cell.
$ws.Cells.Item(24, 18).WrapText = $true
